$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "63.331.97"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +1.29%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.458.35"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.01%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "573.55"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.24%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "146.93"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("E7").Value = "  -0.01%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.539"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.20%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.456.93"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("E10").Value = "  +1.40%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.156"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("E13").Value = "  +0.49%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "27.13"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.18%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.0000181"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.93%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "2.904.27"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.53%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "63.146.64"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +1.11%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.453.87"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.43%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.31"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("E20").Value = "  +5.18%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "329.47"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.68%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.23"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.35%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "2.10"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +14.75%  "
$ws.Range("E24").Value = "  +0.06%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "65.74"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -2.22%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "621.15"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +5.18%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.03"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +5.68%  "
$ws.Range("E28").Value = "  +3.31%  "
$ws.Range("B29").Value = "Fetch.AI"
$ws.Range("C29").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.51"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +4.99%  "
$ws.Range("B30").Value = "WrappedeETH"
$ws.Range("C30").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "2.563.94"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.13%  "
$ws.Range("E31").Value = "  +0.21%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "8.28"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.93%  "
$ws.Range("E33").Value = "  +1.55%  "
$ws.Range("E34").Value = "  -3.26%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.21"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +7.76%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.54"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +2.10%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").Value = "  -0.39%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "5.46"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.07%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "18.91"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("E42").Value = "  -1.63%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.62"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +7.27%  "
$ws.Range("E44").Value = "  -0.14%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "41.84"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.60%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "149.14"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.54%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.78"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.99%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "21.28"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +3.67%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("E50").Value = "  +0.14%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0234"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.18%  "
